$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.174.58'
$ws.Range("E2").Value = '  -0.64%  '

$ws.Range("D3").Value = '1.828.04'
$ws.Range("E3").Value = '  +0.97%  '

$ws.Range("D4").Value = '''0.9995'
$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").Value = '''310.05'
$ws.Range("E5").Value = '  -1.06%  '

$ws.Range("D6").Value = '''0.9981'
$ws.Range("E6").Value = '  -0.24%  '

$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D8").Value = '''0.1008'
$ws.Range("E8").Value = '  +28.05%  '

$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").Value = '''0.3929'
$ws.Range("E9").Value = '  -2.00%  '

$ws.Range("D10").Value = '''1.109'
$ws.Range("E10").Value = '  -0.44%  '

$ws.Range("D11").Value = '''41.49'
$ws.Range("E11").Value = '  +1.35%  '

$ws.Range("D12").Value = '''6.437'
$ws.Range("E12").Value = '  +0.78%  '

$ws.Range("D13").Value = '''20.62'
$ws.Range("E13").Value = '  +1.12%  '

$ws.Range("D14").Value = '''0.9996'
$ws.Range("E14").Value = '  -0.12%  '

$ws.Range("E15").Value = '  +1.36%  '

$ws.Range("D16").Value = '''7.306'
$ws.Range("E16").Value = '  -0.42%  '

$ws.Range("D17").Value = '''0.00001149'
$ws.Range("E17").Value = '  +6.08%  '

$ws.Range("D18").Value = '''92.70'
$ws.Range("E18").Value = '  -0.15%  '

$ws.Range("D19").Value = '''0.06641'
$ws.Range("E19").Value = '  +1.11%  '

$ws.Range("D20").Value = '''0.9979'
$ws.Range("E20").Value = '  -0.20%  '

$ws.Range("D21").Value = '''17.23'
$ws.Range("E21").Value = '  -0.53%  '

$ws.Range("D22").Value = '''6.017'
$ws.Range("E22").Value = '  +0.11%  '

$ws.Range("D23").Value = '28.230.45'
$ws.Range("E23").Value = '  -0.64%  '

$ws.Range("D24").Value = '''11.32'
$ws.Range("E24").Value = '  +1.38%  '

$ws.Range("D25").Value = '''2.238'
$ws.Range("E25").Value = '  +0.05%  '

$ws.Range("D26").Value = '''158.14'
$ws.Range("E26").Value = '  -2.13%  '

$ws.Range("B27").Value = 'EthereumClassic'
$ws.Range("C27").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D27").Value = '''20.80'
$ws.Range("E27").Value = '  +1.27%  '

$ws.Range("B28").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C28").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D28").Value = '2.038.84'
$ws.Range("E28").Value = '  +1.03%  '

$ws.Range("D29").Value = '''2.423'
$ws.Range("E29").Value = '  +0.65%  '

$ws.Range("D30").Value = '''127.35'
$ws.Range("E30").Value = '  -1.18%  '

$ws.Range("D31").Value = '''0.1052'
$ws.Range("E31").Value = '  -3.21%  '

$ws.Range("D32").Value = '''1.041'
$ws.Range("E32").Value = '  -2.57%  '

$ws.Range("D33").Value = '''5.598'
$ws.Range("E33").Value = '  +0.30%  '

$ws.Range("D34").Value = '''3.598'
$ws.Range("E34").Value = '  -1.84%  '

$ws.Range("D35").Value = '''0.06757'
$ws.Range("E35").Value = '  -7.01%  '

$ws.Range("D36").Value = '''9.030'
$ws.Range("E36").Value = '  -1.26%  '

$ws.Range("E37").Value = '  +0.10%  '

$ws.Range("D38").Value = '''0.2148'
$ws.Range("E38").Value = '  -1.32%  '

$ws.Range("E39").Value = '  -1.48%  '

$ws.Range("E40").Value = '  -2.25%  '

$ws.Range("D41").Value = '''0.6221'
$ws.Range("E41").Value = '  +0.41%  '

$ws.Range("D42").Value = '''1.178'
$ws.Range("E42").Value = '  +1.85%  '

$ws.Range("D43").Value = '''0.9982'
$ws.Range("E43").Value = '  -0.19%  '

$ws.Range("E44").Value = '  -0.46%  '

$ws.Range("D45").Value = '''0.5939'
$ws.Range("E45").Value = '  -0.97%  '

$ws.Range("D46").Value = '''3.692'
$ws.Range("E46").Value = '  -1.22%  '

$ws.Range("D47").Value = '''1.268'
$ws.Range("E47").Value = '  -3.52%  '

$ws.Range("D48").Value = '''124.30'
$ws.Range("E48").Value = '  -1.45%  '

$ws.Range("D49").Value = '''1.949'
$ws.Range("E49").Value = '  +0.85%  '

$ws.Range("D50").Value = '''1.182'
$ws.Range("E50").Value = '  -3.24%  '

$ws.Range("D51").Value = '''0.06802'
$ws.Range("E51").Value = '  -0.75%  '
